# Add a new worksheet "2020-11-27" at the end of the workbook, matching the
# layout/style of the other daily attendance sheets, and populate it with
# the day's SpO2/heart-rate readings.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the current last sheet so it lands at the end
# of the tab order (Worksheets.Add defaults to inserting before the active
# sheet otherwise).
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "2020-11-27"

# Header row
$ws.Range("A1").Value = "Sr. No"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Address"
$ws.Range("D1").Value = "Job"
$ws.Range("E1").Value = "Time-Stamp"
$ws.Range("F1").Value = "SpO2_value"
$ws.Range("G1").Value = "Heart-rate"
$ws.Range("H1").Value = "Compensated"
$ws.Range("I1").Value = "Ambient"

# Data rows
$rows = @(
    @{ Row=2; SrNo=1; Name="sachin"; Address="301/Sanskruti-1,Andheri, Mumbai"; Job="Software Engineer"; Time="00:15:25"; SpO2=98.2410200777593;  HR=64.12281009995567; Comp="NA"; Amb="NA" },
    @{ Row=3; SrNo=1; Name="sachin"; Address="301/Sanskruti-1,Andheri, Mumbai"; Job="Software Engineer"; Time="00:18:19"; SpO2=97.95162681117016; HR=0;                 Comp="NA"; Amb="NA" },
    @{ Row=4; SrNo=1; Name="sachin"; Address="301/Sanskruti-1,Andheri, Mumbai"; Job="Software Engineer"; Time="00:20:54"; SpO2=97.87804508522433; HR=66.84032472851912; Comp="NA"; Amb="NA" },
    @{ Row=5; SrNo=1; Name="sachin"; Address="301/Sanskruti-1,Andheri, Mumbai"; Job="Software Engineer"; Time="00:25:04"; SpO2=96.27875888965946; HR=63.44762551053533; Comp="NA"; Amb="NA" },
    @{ Row=6; SrNo=1; Name="sachin"; Address="301/Sanskruti-1,Andheri, Mumbai"; Job="Software Engineer"; Time="00:29:14"; SpO2=97.6024606325315;  HR=66.86134448778704; Comp="NA"; Amb="NA" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.SrNo
    $ws.Cells.Item($row, 2).Value = $r.Name
    $ws.Cells.Item($row, 3).Value = $r.Address
    $ws.Cells.Item($row, 4).Value = $r.Job
    $ws.Cells.Item($row, 5).Value = $r.Time
    $ws.Cells.Item($row, 6).Value = $r.SpO2
    $ws.Cells.Item($row, 7).Value = $r.HR
    $ws.Cells.Item($row, 8).Value = $r.Comp
    $ws.Cells.Item($row, 9).Value = $r.Amb
}

# Styling: header row (A1:I1) and the "Sr. No" column (A2:A6) are bold,
# centered/top-aligned, with a thin box border around each cell - matching
# the style used on every other daily attendance sheet.
$headerRange = $ws.Range("A1:I1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$srNoRange = $ws.Range("A2:A6")
$srNoRange.Font.Bold = $true
$srNoRange.HorizontalAlignment = -4108
$srNoRange.VerticalAlignment = -4160
$srNoRange.Borders.LineStyle = 1
